$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Z2").Value = 11
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 12
$ws.Range("AJ2").Value = 21
$ws.Range("AO2").Value = 8.5
$ws.Range("BB2").Value = 501

# Row 4
$ws.Range("G4").Value = 3.8
$ws.Range("I4").Value = 2.05
$ws.Range("J4").Value = 4.5
$ws.Range("L4").Value = 2.88
$ws.Range("Y4").Value = 15
$ws.Range("AA4").Value = 41
$ws.Range("AC4").Value = 7
$ws.Range("AI4").Value = 8.5
$ws.Range("AJ4").Value = 9.5
$ws.Range("AO4").Value = 23
$ws.Range("AX4").Value = 12
$ws.Range("AY4").Value = 26
$ws.Range("BA4").Value = 67

# Row 6
$ws.Range("G6").Value = 5.5
$ws.Range("I6").Value = 1.7
$ws.Range("J6").Value = 6
$ws.Range("L6").Value = 2.38
$ws.Range("N6").Value = 7.5
$ws.Range("Q6").Value = 2.3
$ws.Range("R6").Value = 1.6
$ws.Range("W6").Value = 12
$ws.Range("X6").Value = 26
$ws.Range("Y6").Value = 19
$ws.Range("AA6").Value = 51
$ws.Range("AD6").Value = 7
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 81
$ws.Range("AK6").Value = 12
$ws.Range("AN6").Value = 7
$ws.Range("AO6").Value = 34
$ws.Range("AQ6").Value = 126
$ws.Range("AW6").Value = 3.5
$ws.Range("AX6").Value = 9

# Row 8
$ws.Range("G8").Value = 2.5
$ws.Range("H8").Value = 3.65
$ws.Range("I8").Value = 2.42
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 2.27
$ws.Range("L8").Value = 2.9
$ws.Range("V8").Value = 2.25
$ws.Range("X8").Value = 14.5
$ws.Range("Z8").Value = 28
$ws.Range("AA8").Value = 18.5
$ws.Range("AD8").Value = 7.4
$ws.Range("AE8").Value = 12
$ws.Range("AJ8").Value = 9.5
$ws.Range("AN8").Value = 4.65
$ws.Range("AO8").Value = 12.5
$ws.Range("AU8").Value = 6.5
$ws.Range("BA8").Value = 70

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("O13").Value = 1.25
